$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "50÷5=10, 0"; New = "79÷5=15, 4" },
    @{ Old = "65÷4=16, 1"; New = "95÷3=31, 2" },
    @{ Old = "69÷7=9, 6";  New = "80÷8=10, 0" },
    @{ Old = "88÷9=9, 7";  New = "31÷5=6, 1" },
    @{ Old = "14÷3=4, 2";  New = "22÷4=5, 2" },
    @{ Old = "58÷7=8, 2";  New = "32÷3=10, 2" },
    @{ Old = "67÷7=9, 4";  New = "23÷2=11, 1" },
    @{ Old = "21÷9=2, 3";  New = "43÷3=14, 1" },
    @{ Old = "56÷6=9, 2";  New = "84÷4=21, 0" },
    @{ Old = "28÷6=4, 4";  New = "29÷4=7, 1" },
    @{ Old = "27÷6=4, 3";  New = "60÷4=15, 0" },
    @{ Old = "60÷6=10, 0"; New = "51÷8=6, 3" },
    @{ Old = "92÷9=10, 2"; New = "93÷7=13, 2" },
    @{ Old = "29÷5=5, 4";  New = "45÷7=6, 3" },
    @{ Old = "44÷6=7, 2";  New = "99÷8=12, 3" },
    @{ Old = "25÷6=4, 1";  New = "43÷5=8, 3" },
    @{ Old = "64÷9=7, 1";  New = "30÷7=4, 2" },
    @{ Old = "50÷4=12, 2"; New = "86÷2=43, 0" },
    @{ Old = "65÷8=8, 1";  New = "59÷9=6, 5" },
    @{ Old = "54÷6=9, 0";  New = "44÷7=6, 2" },
    @{ Old = "37÷4=9, 1";  New = "87÷6=14, 3" },
    @{ Old = "24÷3=8, 0";  New = "81÷3=27, 0" },
    @{ Old = "15÷4=3, 3";  New = "92÷2=46, 0" },
    @{ Old = "35÷7=5, 0";  New = "70÷2=35, 0" },
    @{ Old = "42÷5=8, 2";  New = "16÷5=3, 1" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
